$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarterly indexing bug-fix: column A holds date serials for each row.
# Historically these were set to the 1st day of the quarter-start month;
# the corrected value is the 15th of the following month (mid-quarter
# anchor date). Map of row -> corrected date serial:
$newDates = @{}
$newDates[2] = 25614
$newDates[3] = 25614
$newDates[4] = 25614
$newDates[5] = 25614
$newDates[6] = 25614
$newDates[7] = 25614
$newDates[8] = 25614
$newDates[9] = 25614
$newDates[10] = 25614
$newDates[11] = 25614
$newDates[12] = 39401
$newDates[13] = 39493
$newDates[14] = 39583
$newDates[15] = 39675
$newDates[16] = 39767
$newDates[17] = 39859
$newDates[18] = 39948
$newDates[19] = 40040
$newDates[20] = 40132
$newDates[21] = 40224
$newDates[22] = 40313
$newDates[23] = 40405
$newDates[24] = 40497
$newDates[25] = 40589
$newDates[26] = 40678
$newDates[27] = 40770
$newDates[28] = 40862
$newDates[29] = 40954
$newDates[30] = 41044
$newDates[31] = 41136
$newDates[32] = 41228
$newDates[33] = 41320
$newDates[34] = 41409
$newDates[35] = 41501
$newDates[36] = 41593
$newDates[37] = 41685
$newDates[38] = 41774
$newDates[39] = 41866
$newDates[40] = 41958
$newDates[41] = 42050
$newDates[42] = 42139
$newDates[43] = 42231
$newDates[44] = 42323
$newDates[45] = 42415
$newDates[46] = 42505
$newDates[47] = 42597
$newDates[48] = 42689
$newDates[49] = 42781
$newDates[50] = 42870
$newDates[51] = 42962
$newDates[52] = 43054
$newDates[53] = 43146
$newDates[54] = 43235
$newDates[55] = 43327
$newDates[56] = 43419
$newDates[57] = 43511
$newDates[58] = 43600
$newDates[59] = 43692
$newDates[60] = 43784
$newDates[61] = 43876
$newDates[62] = 43966
$newDates[63] = 44058
$newDates[64] = 44150
$newDates[65] = 44242
$newDates[66] = 44331
$newDates[67] = 44423
$newDates[68] = 44515
$newDates[69] = 44607
$newDates[70] = 44696
$newDates[71] = 44788
$newDates[72] = 44880
$newDates[73] = 44972

foreach ($row in $newDates.Keys) {
    $ws.Cells.Item($row, 1).Value = $newDates[$row]
}
